# The commit swaps the presentation's two theme parts: the theme that
# actually drives the deck (slide master / layouts / slides - serialized
# as ppt/theme/theme2.xml, "Integral" / "Red Violet") takes on the colour
# values that used to live in the otherwise-unused "Office Theme" theme
# part, and vice versa.
#
# The PowerPoint object model only exposes (and only persists writes to)
# the colour scheme of the theme that is actually attached to the slide
# master / slides - there is no supported COM path to the second,
# unattached theme part, and the theme/colour-scheme "Name" properties
# are read-only through this object model. So we reproduce the
# observable, supported part of the change: push the twelve theme colours
# that the "Office Theme" used to hold onto the presentation's live
# colour scheme.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

function ToColorRef([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ColorScheme.Colors(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - i.e. the <a:clrScheme> child order.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $cs.Colors($i).RGB = ToColorRef($officeColors[$i - 1])
}

# Best-effort: try to rename the design/theme to match the swapped
# theme's name too (no-op on hosts where Design.Name is read-only).
try {
    $p.Designs.Item(1).Name = "Office Theme"
} catch {
}

Write-Host "Swapped presentation theme colours to the Office colour scheme."
